$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 719.1316
$ws.Range("I28").Value = 644
$ws.Range("J28").Value = 834.3333
$ws.Range("K28").Value = 644
$ws.Range("L28").Value = 834.3333
$ws.Range("M28").Value = -159
$ws.Range("N28").Value = -1804.3333
$ws.Range("H38").Value = 78.69231000000001
$ws.Range("J38").Value = 100
$ws.Range("L38").Value = 300
$ws.Range("N38").Value = -1044
$ws.Range("H51").Value = 3081.818
$ws.Range("I51").Value = 2500
$ws.Range("J51").Value = 3211.111
$ws.Range("K51").Value = 2500
$ws.Range("L51").Value = 3211.111
$ws.Range("M51").Value = -2016
$ws.Range("N51").Value = -4179.111
$ws.Range("H62").Value = 2323.182
$ws.Range("I62").Value = 2294
$ws.Range("J62").Value = 2385.7144
$ws.Range("K62").Value = 2294
$ws.Range("L62").Value = 2385.7144
$ws.Range("M62").Value = -1670
$ws.Range("N62").Value = -3633.7144
$ws.Range("H65").Value = 2323.182
$ws.Range("I65").Value = 2294
$ws.Range("J65").Value = 2385.7144
$ws.Range("K65").Value = 11470
$ws.Range("L65").Value = 11928.572
$ws.Range("M65").Value = -8350
$ws.Range("N65").Value = -18168.572
$ws.Range("H100").Value = 1524.7826
$ws.Range("I100").Value = 1454.375
$ws.Range("J100").Value = 1685.7142
$ws.Range("K100").Value = 1454.375
$ws.Range("L100").Value = 1685.7142
$ws.Range("M100").Value = -913.375
$ws.Range("N100").Value = -2767.7142
$ws.Range("H116").Value = 2486.4285
$ws.Range("I116").Value = 2376.25
$ws.Range("J116").Value = 2633.3333
$ws.Range("K116").Value = 2376.25
$ws.Range("L116").Value = 2633.3333
$ws.Range("M116").Value = 1065.75
$ws.Range("N116").Value = -9517.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2943031.8
$ws.Range("I2").Value = 2362.8572
$ws.Range("J2").Value = 9804592
$ws.Range("K2").Value = 2362.8572
$ws.Range("L2").Value = 9804592
$ws.Range("M2").Value = -2249.8572
$ws.Range("N2").Value = -9804818
$ws.Range("H45").Value = 1161.9048
$ws.Range("I45").Value = 981.25
$ws.Range("J45").Value = 1273.0769
$ws.Range("K45").Value = 981.25
$ws.Range("L45").Value = 1273.0769
$ws.Range("M45").Value = -604.25
$ws.Range("N45").Value = -2027.0769
$ws.Range("H62").Value = 53000
$ws.Range("J62").Value = 53000
$ws.Range("L62").Value = 53000
$ws.Range("N62").Value = -54248
$ws.Range("H65").Value = 53000
$ws.Range("J65").Value = 53000
$ws.Range("L65").Value = 159000
$ws.Range("N65").Value = -165240
$ws.Range("H97").Value = 812.129
$ws.Range("I97").Value = 715.04
$ws.Range("J97").Value = 1216.6666
$ws.Range("K97").Value = 715.04
$ws.Range("L97").Value = 1216.6666
$ws.Range("M97").Value = -219.04
$ws.Range("N97").Value = -2208.6666
$ws.Range("H110").Value = 1154.4
$ws.Range("I110").Value = 1035.05
$ws.Range("J110").Value = 1631.8
$ws.Range("K110").Value = 1035.05
$ws.Range("L110").Value = 1631.8
$ws.Range("M110").Value = 1009.95
$ws.Range("N110").Value = -5721.8
$ws.Range("H116").Value = 2943031.8
$ws.Range("I116").Value = 2362.8572
$ws.Range("J116").Value = 9804592
$ws.Range("K116").Value = 2362.8572
$ws.Range("L116").Value = 9804592
$ws.Range("M116").Value = -68.85719999999992
$ws.Range("N116").Value = -9809180

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2943031.8
$ws.Range("I3").Value = 2362.8572
$ws.Range("J3").Value = 9804592
$ws.Range("K3").Value = 2362.8572
$ws.Range("L3").Value = 9804592
$ws.Range("M3").Value = -2248.8572
$ws.Range("N3").Value = -9804820
$ws.Range("H105").Value = 2166.6667
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 2250
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 2250
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -5744

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 6500
$ws.Range("J4").Value = 6500
$ws.Range("L4").Value = 6500
$ws.Range("N4").Value = -6724
$ws.Range("H22").Value = 204.11111
$ws.Range("I22").Value = 196.44
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 196.44
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 153.56
$ws.Range("N22").Value = -1000
$ws.Range("H31").Value = 3202.25
$ws.Range("I31").Value = 1558.5883
$ws.Range("J31").Value = 7194
$ws.Range("K31").Value = 1558.5883
$ws.Range("L31").Value = 7194
$ws.Range("M31").Value = -1263.5883
$ws.Range("N31").Value = -7784
$ws.Range("H34").Value = 3202.25
$ws.Range("I34").Value = 1558.5883
$ws.Range("J34").Value = 7194
$ws.Range("K34").Value = 1558.5883
$ws.Range("L34").Value = 7194
$ws.Range("M34").Value = -1356.5883
$ws.Range("N34").Value = -7598
$ws.Range("H107").Value = 4168185.5
$ws.Range("I107").Value = 8335091
$ws.Range("K107").Value = 8335091
$ws.Range("M107").Value = -8333171
$ws.Range("H122").Value = 11230.2
$ws.Range("I122").Value = 12255.777
$ws.Range("K122").Value = 36767.331
$ws.Range("M122").Value = -34317.331

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 12079.8
$ws.Range("I39").Value = 2999
$ws.Range("J39").Value = 14350
$ws.Range("K39").Value = 8997
$ws.Range("L39").Value = 43050
$ws.Range("M39").Value = -8703
$ws.Range("N39").Value = -43638
$ws.Range("H110").Value = 1339.7142
$ws.Range("I110").Value = 1339.7142
$ws.Range("K110").Value = 4019.1426
$ws.Range("M110").Value = 70.85740000000033
$ws.Range("H120").Value = 16059.5
$ws.Range("I120").Value = 7738
$ws.Range("J120").Value = 18833.334
$ws.Range("K120").Value = 23214
$ws.Range("L120").Value = 56500.00199999999
$ws.Range("M120").Value = -18376
$ws.Range("N120").Value = -66176.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1919.4231
$ws.Range("J5").Value = 1931.3636
$ws.Range("L5").Value = 1931.3636
$ws.Range("N5").Value = -2155.3636
$ws.Range("H102").Value = 2571.4285
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 2800
$ws.Range("K102").Value = 2000
$ws.Range("L102").Value = 2800
$ws.Range("M102").Value = -378
$ws.Range("N102").Value = -6044
$ws.Range("H113").Value = 3651
$ws.Range("I113").Value = 10980
$ws.Range("J113").Value = 2185.2
$ws.Range("K113").Value = 10980
$ws.Range("L113").Value = 2185.2
$ws.Range("M113").Value = -8810
$ws.Range("N113").Value = -6525.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 100000
$ws.Range("J62").Value = 100000
$ws.Range("L62").Value = 100000
$ws.Range("N62").Value = -101248
$ws.Range("H65").Value = 100000
$ws.Range("J65").Value = 100000
$ws.Range("L65").Value = 300000
$ws.Range("N65").Value = -306240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1770
$ws.Range("I122").Value = 1744.4445
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5233.333500000001
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2783.333500000001
$ws.Range("N122").Value = -10900
